$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new client code (40125) to the RECURRENTE_100K-200K group (row 6, column B)
$currentValue = $ws.Range("B6").Value2
$ws.Range("B6").Value = $currentValue + ".40125"

# Update the selected/active cell from B10 to B7
$ws.Range("B7").Select()
